# LOQ4249.xlsx edit: re-sync the "Oficina de Inovacao" course-description
# sheet. Rows 10-22 are rebuilt from scratch (upstream shared-string churn
# caused row 22 to disappear and rows 13-21 to shift/change content).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- remove the old rows 10..22 completely (bottom-up so indices stay valid)
for ($r = 22; $r -ge 10; $r--) {
    $ws.Rows.Item($r).Delete()
}

# --- helper: write a whole data row (A/B/C) and fix up height + column-B style
function Set-DataRow($Row, $A, $B, $C, $Height) {

    if ($A -ne $null) { $ws.Range("A$Row").Value = $A }

    if ($B -ne $null) {
        # Column B's <col> definitions overlap (A:B style 1, then B alone
        # style 2), so a brand-new B cell inherits the wrong style unless we
        # copy correct formatting over from an existing, known-good B cell.
        $ws.Range("B2").Copy() | Out-Null
        $ws.Range("B$Row").PasteSpecial(-4122) | Out-Null
        $ws.Range("B$Row").Value = $B
    }

    if ($C -ne $null) { $ws.Range("C$Row").Value = $C }

    if ($Height -gt 0) {
        $ws.Rows.Item($Row).RowHeight = $Height
    }
}

Set-DataRow 10 "Objetivos:" "5840560 - Marco Antonio Carvalho Pereira" "5840560 - Marco Antonio Carvalho Pereira" 60
Set-DataRow 11 "Objectives:" $null $null 60
Set-DataRow 12 "Docentes responsáveis:" $null $null 0
Set-DataRow 13 "Programa resumido:" "Semestral" "Semestral" 60
Set-DataRow 14 "Short syllabus:" $null $null 60
Set-DataRow 15 "Programa:" "01/01/2018" "01/01/2018" 120
Set-DataRow 16 "Syllabus:" $null $null 120
Set-DataRow 17 "Avaliação:" $null $null 0
Set-DataRow 18 "Método:" "5840560 - Marco Antonio Carvalho Pereira" "5840560 - Marco Antonio Carvalho Pereira" 60
Set-DataRow 19 "Critério:" "Atividades docentes: Mentoria, palestras e seminários.Atividades discentes: Elaboração de projeto utilizando laboratórios e instalações da USP." "Atividades docentes: Mentoria, palestras e seminários.Atividades discentes: Elaboração de projeto utilizando laboratórios e instalações da USP." 60
Set-DataRow 20 "Norma de recuperação:" "Avaliação pela equipe de mentores, considerando critérios, tais como: qualidade técnica da proposta, grau de inovação, viabilidade técnica, dentre outros.Nota de projeto maior ou igual a 5,0 (cinco)." "Avaliação pela equipe de mentores, considerando critérios, tais como: qualidade técnica da proposta, grau de inovação, viabilidade técnica, dentre outros.Nota de projeto maior ou igual a 5,0 (cinco)." 60
Set-DataRow 21 "Bibliografia:" "Não há recuperação." "Não há recuperação." 120
